$d = $word.ActiveDocument

# The placeholder "<<ThiSinh_DienThoai>>" appears twice in the document:
#   - paragraph "Dia chi lien lac khi can bao tin: <<ThiSinh_DienThoai>>"
#   - paragraph "Dien thoai lien lac (de Truong bao tin): <<ThiSinh_DienThoai>>"
# Only the first one (the "dia chi" / address line) should be updated to
# "<<ThiSinh_DCNhanGiayBao>>", with the run now carrying an explicit
# black (000000) font color. The second occurrence must stay untouched.

$para = $d.Paragraphs.Item(14)
$range = $para.Range

$range.Find.ClearFormatting()
$range.Find.Replacement.ClearFormatting()
$range.Find.Replacement.Font.Color = 0

$range.Find.Execute("ThiSinh_DienThoai", $true, $false, $false, $false, $false, `
                     $true, 1, $false, "ThiSinh_DCNhanGiayBao", 2)
